$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns of interest: D=4 (Fecha), J=10 (Volumen), K=11 (Precio minimo),
# L=12 (Precio maximo), M=13 (Precio promedio ponderado), P=16 (Precio $/Kg)
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot current ("before") values for every affected source row, per column,
# before any writes happen. Use Value2 (plain property) instead of Value
# (parameterized property) so the underlying primitive is returned/stored correctly.
$snap = @{}
$srcRows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 22, 23, 24, 25, 26, 27)
foreach ($r in $srcRows) {
    foreach ($c in $cols) {
        $snap["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Row permutation: new row (key) gets the old values of the row (value).
$mapping = @{
    3 = 15
    4 = 23
    5 = 20
    6 = 7
    7 = 16
    8 = 17
    9 = 5
    10 = 19
    11 = 13
    12 = 11
    13 = 27
    14 = 25
    15 = 18
    16 = 4
    17 = 24
    18 = 14
    19 = 3
    20 = 6
    22 = 8
    23 = 12
    24 = 26
    25 = 10
    26 = 22
    27 = 9
}

foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    foreach ($c in $cols) {
        $ws.Cells.Item($dst, $c).Value2 = $snap["$src-$c"]
    }
}

$wb.Save()
